# Revert: re-add the "TestProject" / "Test" row that was removed.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row is appended right after the current last row (row 4) -> row 5
$row = 5

$ws.Cells.Item($row, 1).Value = 3
$ws.Cells.Item($row, 2).Value = "TestProject"
$ws.Cells.Item($row, 3).Value = "Test"

$ws.Cells.Item($row, 4).Value = 45772
$ws.Cells.Item($row, 4).NumberFormat = "dd/MM/yyyy"

$ws.Cells.Item($row, 5).Value = 45773
$ws.Cells.Item($row, 5).NumberFormat = "dd/MM/yyyy"

$ws.Cells.Item($row, 6).Value = 5
$ws.Cells.Item($row, 7).Value = $false
